$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Expand the existing table with a new "Authorship Resource" column (11th column)
$table = $ws.ListObjects.Item("Table1")
$table.ListColumns.Add() | Out-Null

# Header + data for the new column K
$ws.Range("K1").Value = "Authorship Resource"
$ws.Range("K2").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("K3").Value = "Daniela Subotic, Noémi Villars-Amberg"

# Copy the formatting from the neighbouring "Book Cover ID" column so the
# new header/data cells render the same way (wrap text, top-aligned, etc.)
$ws.Range("J1:J3").Copy() | Out-Null
$ws.Range("K1:K3").PasteSpecial(-4122) | Out-Null

# Reflect where the user last clicked after adding the column
$ws.Range("L11").Select() | Out-Null
